$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# --- Red deer instantaneous death-rate / birth-rate columns (S:W), rows 15-35 ---
# Labels first typed at row 15 (top of table) and row 40 (summary), then the
# per-year ratio formulas are filled down alongside.
$ws.Range("S15").Value = "death red deer:"
$ws.Range("S40").Value = "instantaneous dr"
$ws.Range("V15").Value = "birth red deer:"
$ws.Range("V40").Value = "instantaneous birth rate:"
$ws.Range("F42").Value = "red deer :"
$ws.Range("T43").Value = "r = "
$ws.Range("T37").Value = "don""t count the last years bc of the capacity and something obviously happened there."

for ($r = 15; $r -le 35; $r++) {
    $ws.Range("T$r").Formula = "=J$r/D$r"
    $ws.Range("W$r").Formula = "=P$r/D$r"
}

# Blank, but styled, spacer cells below the filled-down ratio columns
# (copy the bold "Total mortality" header formatting from H3).
$ws.Range("H3").Copy()
$ws.Range("T39").PasteSpecial(-4122)
$ws.Range("W39").PasteSpecial(-4122)

# Averages of the instantaneous rates (years 1992-2012 only; the comment in
# T37 explains why the last few years are excluded).
$ws.Range("T40").Formula = "=AVERAGE(T15:T38)"
$ws.Range("W40").Formula = "=AVERAGE(W15:W38) + 1"

# Difference between the (offset) birth rate and death rate -> intrinsic growth rate r.
$ws.Range("U43").Formula = "=W40-T40"

# Leave the selection where the author left it when they saved the file.
$null = $ws.Range("T37").Select()
